$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- workbook.xml: add new (non-hidden) "_xlnm._FilterDatabase_0" defined names,
#     one per sheet, mirroring the existing visible "_xlnm._FilterDatabase" entries ---
$ws1.Names.Add("_xlnm._FilterDatabase_0", "='Main root'!`$A`$1:`$B`$2")
$ws2.Names.Add("_xlnm._FilterDatabase_0", "=Tests!`$A`$1:`$C`$4")

# --- xl/sharedStrings.xml: introduce the new "Identifier" string. It lands on
#     "Main root"!A1 (previously "Id"), which is the one cell that ends up
#     pointing at it once the shared-string table gains a new entry ---
$ws1.Range("A1").Value = "Identifier"

# --- xl/worksheets/sheet2.xml: move the "Tests" sheet's bottom-right pane
#     selection from B12 to A2 ---
[void]$ws2.Range("A2").Select()

# Restore "Main root" as the active tab (selecting on Tests switches focus there)
[void]$ws1.Activate()
